# Weekly update: insert a new price observation as row 107, pushing the
# existing rows 107-136 down to 108-137 (new row 137 is therefore a
# duplicate of what used to be the last row, 136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 107; this shifts rows 107-136
# down to 108-137 and copies formatting (incl. the date style) from the
# row above.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A107").Value = 11
$ws.Range("B107").Value = "Vega Monumental Concepción"
$ws.Range("C107").Value = "Bíobío"
$ws.Range("D107").Value = "2022-06-24"
$ws.Range("E107").Value = 8
$ws.Range("F107").Value = 100112043
$ws.Range("G107").Value = "Pepino ensalada"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 18000
$ws.Range("L107").Value = 20000
$ws.Range("M107").Value = 19000
$ws.Range("N107").Value = "`$/caja 60 unidades"
$ws.Range("O107").Value = "Región de Arica y Parinacota"
$ws.Range("P107").Value = 317
$ws.Range("Q107").Value = 60
$ws.Range("R107").Value = "Hortaliza"
